# feat: now it generates vehicle columns
#
# Applies:
#  - Entrada!C1 "P001" -> "1" (kept as text, same style)
#  - Relatório: reshuffle the "Data"/"Movimentos" summary up one row,
#    rename labels, change values, and add a new vehicle-count header
#    table spanning columns B:AD across rows 3-4 (merged group headers
#    in row 3, sub-headers in row 4), with matching column widths,
#    borders and centered alignment.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Entrada": Ponto value P001 -> 1 (stays text, keeps style)
# ---------------------------------------------------------------
$wsEntrada = $wb.Worksheets.Item("Entrada")

$tmp = $wsEntrada.Range("Z1")
$tmp.NumberFormat = "@"
$tmp.Value2 = "1"
$tmp.Copy()
$wsEntrada.Range("C1").PasteSpecial(-4163)  # xlPasteValues
$tmp.Clear()

# ---------------------------------------------------------------
# Sheet "Relatório": rebuild the summary + vehicle table
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Relatório")

# --- Row 1: "Data:" / date value (keep old C2's style, keep text literal) ---
$ws.Range("C2").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

$tmp2 = $ws.Range("AZ1")
$tmp2.NumberFormat = "@"
$tmp2.Value2 = "04-05-2025"
$tmp2.Copy()
$ws.Range("C1").PasteSpecial(-4163)  # xlPasteValues
$tmp2.Clear()

$ws.Range("B1").Value2 = "Data:"

# --- Row 2: "Movimento:" / "1A" (reuses existing C2 style/content cell) ---
$ws.Range("B2").Value2 = "Movimento:"
$ws.Range("C2").Value2 = "1A"

# --- Clear old row 3 ("Movimentos" / "P001A") before rebuilding ---
$ws.Range("B3:C3").Clear()

# --- Column widths (characters); COM pads +5/6, so subtract that back out ---
$offset = 0.8333333333333334
$colWidths = @{
    4  = 7
    5  = 12
    6  = 9
    7  = 9
    8  = 5
    9  = 11
    10 = 9
    11 = 9
    12 = 9
    13 = 5
    14 = 5
    15 = 5
    16 = 5
    17 = 5
    18 = 5
    19 = 8
    20 = 10
    21 = 7
    22 = 9
    23 = 11
    24 = 8
    25 = 10
    26 = 8
    27 = 8
    28 = 7
    29 = 7
    30 = 17
}
foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col] - $offset
}

# --- Row 3 formatting: medium box border, centered, default font ---
$row3 = $ws.Range("B3:AD3")
$row3.HorizontalAlignment = -4108  # xlCenter
$row3.VerticalAlignment = -4108    # xlCenter
$row3.Borders.LineStyle = 1
$row3.Borders.Weight = -4138       # xlMedium

# --- Row 4 "spanning" placeholders (share row-3 style, vertically merged) ---
$row4Thick = $ws.Range("D4,H4,U4,AD4")
foreach ($a in $row4Thick.Areas) {
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4108
    $a.Borders.LineStyle = 1
    $a.Borders.Weight = -4138
}

# --- Row 4 sub-headers: thin box border, centered, 10pt font ---
$row4Thin = $ws.Range("B4:C4,E4:G4,I4:T4,V4:AC4")
foreach ($a in $row4Thin.Areas) {
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4108
    $a.Borders.LineStyle = 1
    $a.Borders.Weight = 2          # xlThin
    $a.Font.Size = 10
}

# --- Row 3 header text ---
$ws.Range("B3").Value2 = "Horas"
$ws.Range("D3").Value2 = "Leves"
$ws.Range("E3").Value2 = "Carretinha"
$ws.Range("H3").Value2 = "VUC"
$ws.Range("I3").Value2 = "Caminhões"
$ws.Range("L3").Value2 = "Carreta"
$ws.Range("S3").Value2 = "Ônibus"
$ws.Range("U3").Value2 = "Motos"
$ws.Range("V3").Value2 = "Pesados"
$ws.Range("AD3").Value2 = "Veículos Totais"

# --- Row 4 sub-header text ---
$ws.Range("B4").Value2 = "das"
$ws.Range("C4").Value2 = "as"
$ws.Range("E4").Value2 = "1 Eixo"
$ws.Range("F4").Value2 = "2 Eixos"
$ws.Range("G4").Value2 = "3 Eixos"
$ws.Range("I4").Value2 = "2 Eixos"
$ws.Range("J4").Value2 = "3 Eixos"
$ws.Range("K4").Value2 = "4 Eixos"
$ws.Range("L4").Value2 = "2 E"
$ws.Range("M4").Value2 = "3 E"
$ws.Range("N4").Value2 = "4 E"
$ws.Range("O4").Value2 = "5 E"
$ws.Range("P4").Value2 = "6 E"
$ws.Range("Q4").Value2 = "7 E"
$ws.Range("R4").Value2 = "8 E"
$ws.Range("S4").Value2 = "2 E"
$ws.Range("T4").Value2 = "3 E ou +"
$ws.Range("V4").Value2 = "% Cam"
$ws.Range("W4").Value2 = "Caminhões"
$ws.Range("X4").Value2 = "% Carr"
$ws.Range("Y4").Value2 = "Carretas"
$ws.Range("Z4").Value2 = "% Ônib"
$ws.Range("AA4").Value2 = "Ônibus"
$ws.Range("AB4").Value2 = "% Pes"
$ws.Range("AC4").Value2 = "Total"

# --- Merges ---
$ws.Range("B3:C3").Merge()
$ws.Range("D3:D4").Merge()
$ws.Range("E3:G3").Merge()
$ws.Range("H3:H4").Merge()
$ws.Range("I3:K3").Merge()
$ws.Range("L3:R3").Merge()
$ws.Range("S3:T3").Merge()
$ws.Range("U3:U4").Merge()
$ws.Range("V3:AC3").Merge()
$ws.Range("AD3:AD4").Merge()
